# Insert a new row of weekly price data at row 6, pushing existing
# rows 6-12 down to 7-13, then populate the new row 6 with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6 (shifts rows 6..12 -> 7..13)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the newly reported week's values.
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 45037
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 100112041
$ws.Cells.Item(6, 7).Value = "Fruto del paraíso"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(6, 11).Value = 24000
$ws.Cells.Item(6, 12).Value = 24000
$ws.Cells.Item(6, 13).Value = 24000
$ws.Cells.Item(6, 14).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 1600
$ws.Cells.Item(6, 17).Value = 15
$ws.Cells.Item(6, 18).Value = "Hortaliza"
